$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    # Assigning a date-like literal string (e.g. "2026-02-17") straight to
    # .Value lets the COM layer auto-coerce it into a date serial, same as
    # typing it into a pre-formatted General cell in real Excel. Forcing
    # the cell to Text first keeps it a literal string; ClearFormats()
    # afterwards drops the now-unneeded explicit number format so the
    # cell is left with no special style, matching a plain literal cell.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

function Set-EmptyCell($ws, $row, $col) {
    # A plain "" assignment to a previously-untouched cell leaves no trace
    # in the saved sheet at all (the cell is dropped rather than kept as
    # an empty-but-present record). Routing it through the Text-format
    # dance keeps the cell materialised (an empty <c> element) the same
    # way the source row stores its blank "Exit Price" / "Exit Reason"
    # cells.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = ""
    $c.ClearFormats()
}

# ---------------------------------------------------------------------
# "All Trades" sheet: append trade #65 as a new row 66
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Cells.Item(66, 1).Value = 65
Set-TextCell $wsAll 66 2 "2026-02-17"
Set-TextCell $wsAll 66 3 "20:47:48"
Set-TextCell $wsAll 66 4 "MarketMaking"
Set-TextCell $wsAll 66 5 "UP"
$wsAll.Cells.Item(66, 6).Value = 0.16
Set-EmptyCell $wsAll 66 7
Set-TextCell $wsAll 66 8 "OPEN"
$wsAll.Cells.Item(66, 9).Value = 0
$wsAll.Cells.Item(66, 10).Value = 0
$wsAll.Cells.Item(66, 11).Value = 100.32
Set-EmptyCell $wsAll 66 12
$wsAll.Cells.Item(66, 13).Value = 0
$wsAll.Cells.Item(66, 14).Value = 0
$wsAll.Cells.Item(66, 15).Value = 0
$wsAll.Cells.Item(66, 16).Value = 0.6
Set-TextCell $wsAll 66 17 "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# "MarketMaking" sheet: append the same trade as new row 33
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

$wsMM.Cells.Item(33, 1).Value = 65
Set-TextCell $wsMM 33 2 "2026-02-17"
Set-TextCell $wsMM 33 3 "20:47:48"
Set-TextCell $wsMM 33 4 "MarketMaking"
Set-TextCell $wsMM 33 5 "UP"
$wsMM.Cells.Item(33, 6).Value = 0.16
Set-EmptyCell $wsMM 33 7
Set-TextCell $wsMM 33 8 "OPEN"
$wsMM.Cells.Item(33, 9).Value = 0
$wsMM.Cells.Item(33, 10).Value = 0
$wsMM.Cells.Item(33, 11).Value = 100.32
$wsMM.Cells.Item(33, 12).Value = 0
$wsMM.Cells.Item(33, 13).Value = 0
$wsMM.Cells.Item(33, 14).Value = 0.6
Set-TextCell $wsMM 33 15 "Normal spread capture: 19600 bps"
Set-EmptyCell $wsMM 33 16
$wsMM.Cells.Item(33, 17).Value = 0
